# SCRUM Update E, G, H und I
#
# Applies the changes described by the commit:
#  - Backlog sheet (sheet1): extend several "Beschreibung" (description)
#    texts in column C with more detail.
#  - Sprint-Backlog sheet (sheet2): add a new column E ("Done"/"10min" style
#    status marker), fill in Due-date / status for rows 3-4, and populate
#    rows 5-7 with the next three backlog items (estimation + status),
#    matching the longer descriptions used on the Backlog sheet.
#  - Update the remembered cursor/selection position on both sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Backlog"
$ws2 = $wb.Worksheets.Item(2)   # "Sprint-Backlog"

# ---------------------------------------------------------------------
# Sheet "Backlog": lengthen the description texts in column C (rows 6-9)
# ---------------------------------------------------------------------
$ws1.Range("C6").Value = "Umsetzung des Painters im Framework. Standard TicTacToe. Kreise und Kreuze. Ausblick: Animationen."
$ws1.Range("C7").Value = "Umsetzung der Rules im Framework. Standard TicTacToe. Ausblick: Zeitlimit."
$ws1.Range("C8").Value = "Umsetzung eines menschlichen Players im Framework. Satndard TicTacToe mit Mauseingabe"
$ws1.Range("C9").Value = "Umsetzung eines PC gesteuerten Players im Framework. Standard ticTacToe als Ausblick."

# ---------------------------------------------------------------------
# Sheet "Sprint-Backlog": new column E + filled-in rows 3-7
# ---------------------------------------------------------------------

# New column E width (matches the other default-width columns) and its
# centred (no wrap) default alignment.
$ws2.Columns("E").ColumnWidth = 8.88671875

# Header row: D1 gets the "Focus Faktor" note.
$ws2.Range("D1").Value = "Focus Faktor : 0,5 "

# Rows 3 & 4: fill in the due date (already present) + new "Done" status
# in the new column E.
$ws2.Range("D3").Value = "Due: 08.11.2021"
$ws2.Range("E3").Value = "Done"
$ws2.Range("E3").HorizontalAlignment = -4108  # xlCenter
$ws2.Range("E3").VerticalAlignment = -4108    # xlCenter

$ws2.Range("D4").Value = "Due: 08.11.2021"
$ws2.Range("E4").Value = "Done"
$ws2.Range("E4").HorizontalAlignment = -4108
$ws2.Range("E4").VerticalAlignment = -4108

# Row 5: "Spezifikation einer Umsetzungsidee..." item, estimated 10min, Done.
$ws2.Range("A5").Value = 2
$ws2.Range("B5").Value = "Spezifikation einer Umsetzungsidee für das Spiel TicTacToe"
$ws2.Range("C5").Value = "Spezifikation: Wie soll das TicTacToe aussehen? Wie soll das Aussehen technisch erreicht werden? Wie sollen Regeln umgesetzt werden? Soll es Sonderregeln geben, bzw. Sonderspielfelder, etc.? Welche Spieler (PC, ...) soll es geben? Etc.?"
$ws2.Range("D5").Value = "10min"
$ws2.Range("E5").Value = "Done"
$ws2.Range("E5").HorizontalAlignment = -4108
$ws2.Range("E5").VerticalAlignment = -4108

# Row 6: "Implementierung der TicTacToe-Darstellung" item, 240min estimate.
$ws2.Range("A6").Value = 2
$ws2.Range("B6").Value = "Implementierung der TicTacToe-Darstellung"
$ws2.Range("C6").Value = "Umsetzung des Painters im Framework. Standard TicTacToe. Kreise und Kreuze. Ausblick: Animationen."
$ws2.Range("D6").Value = "240min"

# Row 7: "Implementierung eines menschlichen TicTacToe-Spielers" item, 120min.
$ws2.Range("A7").Value = 2
$ws2.Range("B7").Value = "Implementierung eines menschlichen TicTacToe-Spielers"
$ws2.Range("C7").Value = "Umsetzung eines menschlichen Players im Framework. Satndard TicTacToe mit Mauseingabe"
$ws2.Range("D7").Value = "120min"

# ---------------------------------------------------------------------
# Selection / scroll position bookkeeping (matches the saved view state)
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A8:XFD8").Select()
$excel.ActiveWindow.ScrollRow = 2

$ws2.Activate()
$ws2.Range("F3").Select()
